{"js": "// Replace the worksheet date and every \"a\u00f7b=\" problem with the new values\n// from the day's refreshed set. Each \"old\" string is unique in the\n// document, so a scoped, case-sensitive search-and-replace per pair is\n// unambiguous. Pairs are applied in the same top-to-bottom order they\n// appear in the document so that a later pair's \"new\" text (e.g.\n// \"43\u00f73=\") can never be accidentally re-matched by an earlier pair whose\n// \"old\" text is the same string (e.g. \"43\u00f73=\" -> \"83\u00f72=\" runs before\n// \"18\u00f77=\" -> \"43\u00f73=\").\nconst replacements = [\n  [\"2024-02-15 Thursday\", \"2024-02-16 Friday\"],\n  [\"80\u00f73=\", \"94\u00f79=\"],\n  [\"10\u00f74=\", \"22\u00f78=\"],\n  [\"50\u00f78=\", \"77\u00f79=\"],\n  [\"26\u00f73=\", \"71\u00f79=\"],\n  [\"89\u00f75=\", \"70\u00f73=\"],\n  [\"77\u00f72=\", \"37\u00f78=\"],\n  [\"19\u00f72=\", \"49\u00f76=\"],\n  [\"93\u00f78=\", \"55\u00f73=\"],\n  [\"84\u00f76=\", \"57\u00f73=\"],\n  [\"89\u00f72=\", \"99\u00f77=\"],\n  [\"96\u00f73=\", \"86\u00f76=\"],\n  [\"43\u00f73=\", \"83\u00f72=\"],\n  [\"52\u00f77=\", \"45\u00f72=\"],\n  [\"53\u00f75=\", \"60\u00f79=\"],\n  [\"11\u00f73=\", \"77\u00f78=\"],\n  [\"18\u00f77=\", \"43\u00f73=\"],\n  [\"56\u00f75=\", \"77\u00f74=\"],\n  [\"48\u00f76=\", \"88\u00f76=\"],\n  [\"39\u00f75=\", \"94\u00f74=\"],\n  [\"87\u00f79=\", \"64\u00f78=\"],\n  [\"20\u00f77=\", \"55\u00f72=\"],\n  [\"54\u00f75=\", \"21\u00f78=\"],\n  [\"77\u00f77=\", \"49\u00f73=\"],\n  [\"65\u00f76=\", \"55\u00f72=\"],\n  [\"26\u00f79=\", \"82\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Refresh the worksheet date and every \"a\u00f7b=\" division problem with the\n# day's new values. Each \"old\" string occurs exactly once in the document,\n# so Find/Replace (wdReplaceOne) per pair is unambiguous. The pairs are\n# applied in the same top-to-bottom order they appear in the document so a\n# later pair's replacement text (e.g. \"43\u00f73=\") is never re-matched by an\n# earlier pair that searches for that same string (e.g. \"43\u00f73=\" -> \"83\u00f72=\"\n# runs before \"18\u00f77=\" -> \"43\u00f73=\").\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-15 Thursday\", \"2024-02-16 Friday\"),\n    @(\"80\u00f73=\", \"94\u00f79=\"),\n    @(\"10\u00f74=\", \"22\u00f78=\"),\n    @(\"50\u00f78=\", \"77\u00f79=\"),\n    @(\"26\u00f73=\", \"71\u00f79=\"),\n    @(\"89\u00f75=\", \"70\u00f73=\"),\n    @(\"77\u00f72=\", \"37\u00f78=\"),\n    @(\"19\u00f72=\", \"49\u00f76=\"),\n    @(\"93\u00f78=\", \"55\u00f73=\"),\n    @(\"84\u00f76=\", \"57\u00f73=\"),\n    @(\"89\u00f72=\", \"99\u00f77=\"),\n    @(\"96\u00f73=\", \"86\u00f76=\"),\n    @(\"43\u00f73=\", \"83\u00f72=\"),\n    @(\"52\u00f77=\", \"45\u00f72=\"),\n    @(\"53\u00f75=\", \"60\u00f79=\"),\n    @(\"11\u00f73=\", \"77\u00f78=\"),\n    @(\"18\u00f77=\", \"43\u00f73=\"),\n    @(\"56\u00f75=\", \"77\u00f74=\"),\n    @(\"48\u00f76=\", \"88\u00f76=\"),\n    @(\"39\u00f75=\", \"94\u00f74=\"),\n    @(\"87\u00f79=\", \"64\u00f78=\"),\n    @(\"20\u00f77=\", \"55\u00f72=\"),\n    @(\"54\u00f75=\", \"21\u00f78=\"),\n    @(\"77\u00f77=\", \"49\u00f73=\"),\n    @(\"65\u00f76=\", \"55\u00f72=\"),\n    @(\"26\u00f79=\", \"82\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n"}
